# Actualización Automática de Datos (EA1, EA2 y EA3)
# Update the "timestamp" column (H) for all data rows to the new refresh time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-03-26 21:07:25"

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 8).Value = $newTimestamp
}
